$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.173134328358209
$ws.Range("C2").Value = 0.6059701492537314
$ws.Range("J2").Value = 0.008955223880597015
$ws.Range("P2").Value = 0.1492537313432836
$ws.Range("S2").Value = 0.0626865671641791
$ws.Range("C3").Value = 0.04205607476635514
$ws.Range("J3").Value = 0.03271028037383177
$ws.Range("P3").Value = 0.7523364485981309
$ws.Range("S3").Value = 0.1728971962616822
$ws.Range("J4").Value = 0.04615384615384616
$ws.Range("P4").Value = 0.7692307692307693
$ws.Range("S4").Value = 0.1846153846153846
$ws.Range("B6").Value = 0.08203125
$ws.Range("D6").Value = 0.015625
$ws.Range("F6").Value = 0.07421875
$ws.Range("J6").Value = 0.203125
$ws.Range("O6").Value = 0.00390625
$ws.Range("Q6").Value = 0.1796875
$ws.Range("R6").Value = 0.10546875
$ws.Range("S6").Value = 0.3359375
$ws.Range("B7").Value = 0.1030534351145038
$ws.Range("D7").Value = 0.04198473282442748
$ws.Range("F7").Value = 0.05343511450381679
$ws.Range("J7").Value = 0.1106870229007634
$ws.Range("O7").Value = 0.007633587786259542
$ws.Range("Q7").Value = 0.1679389312977099
$ws.Range("R7").Value = 0.1068702290076336
$ws.Range("S7").Value = 0.4083969465648855
$ws.Range("B8").Value = 0.06920415224913495
$ws.Range("D8").Value = 0.02422145328719723
$ws.Range("F8").Value = 0.04844290657439446
$ws.Range("J8").Value = 0.09688581314878893
$ws.Range("O8").Value = 0.01903114186851211
$ws.Range("Q8").Value = 0.198961937716263
$ws.Range("R8").Value = 0.1332179930795848
$ws.Range("S8").Value = 0.4100346020761246
$ws.Range("B9").Value = 0.04597701149425287
$ws.Range("D9").Value = 0.01532567049808429
$ws.Range("E9").Value = 0.003831417624521073
$ws.Range("F9").Value = 0.04597701149425287
$ws.Range("J9").Value = 0.1149425287356322
$ws.Range("O9").Value = 0.01915708812260536
$ws.Range("Q9").Value = 0.2030651340996169
$ws.Range("R9").Value = 0.1187739463601533
$ws.Range("S9").Value = 0.4329501915708812
$ws.Range("B10").Value = 0.1047381546134663
$ws.Range("D10").Value = 0.02057356608478803
$ws.Range("E10").Value = 0.001246882793017456
$ws.Range("F10").Value = 0.06920199501246883
$ws.Range("J10").Value = 0.1097256857855362
$ws.Range("O10").Value = 0.01496259351620948
$ws.Range("Q10").Value = 0.2437655860349127
$ws.Range("R10").Value = 0.09226932668329177
$ws.Range("S10").Value = 0.3435162094763092
$ws.Range("G11").Value = 0.1491442542787286
$ws.Range("J11").Value = 0.09535452322738386
$ws.Range("K11").Value = 0.2102689486552567
$ws.Range("L11").Value = 0.5330073349633252
$ws.Range("S11").Value = 0.01222493887530562
$ws.Range("G12").Value = 0.7324561403508771
$ws.Range("J12").Value = 0.2017543859649123
$ws.Range("K12").Value = 0.01754385964912281
$ws.Range("L12").Value = 0.03070175438596491
$ws.Range("S12").Value = 0.01754385964912281
$ws.Range("G13").Value = 0.6885245901639344
$ws.Range("J13").Value = 0.2295081967213115
$ws.Range("S13").Value = 0.08196721311475409
$ws.Range("F15").Value = 0.02262443438914027
$ws.Range("H15").Value = 0.1583710407239819
$ws.Range("I15").Value = 0.04977375565610859
$ws.Range("J15").Value = 0.3665158371040724
$ws.Range("K15").Value = 0.08144796380090498
$ws.Range("M15").Value = 0.01809954751131222
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.2398190045248869
$ws.Range("F16").Value = 0.007874015748031496
$ws.Range("H16").Value = 0.1929133858267716
$ws.Range("I16").Value = 0.05905511811023622
$ws.Range("J16").Value = 0.405511811023622
$ws.Range("K16").Value = 0.1377952755905512
$ws.Range("M16").Value = 0.02755905511811024
$ws.Range("O16").Value = 0.03149606299212598
$ws.Range("S16").Value = 0.1377952755905512
$ws.Range("F17").Value = 0.01395348837209302
$ws.Range("H17").Value = 0.1736434108527132
$ws.Range("I17").Value = 0.1116279069767442
$ws.Range("J17").Value = 0.4124031007751938
$ws.Range("K17").Value = 0.09147286821705426
$ws.Range("M17").Value = 0.02325581395348837
$ws.Range("O17").Value = 0.04496124031007752
$ws.Range("S17").Value = 0.1286821705426357
$ws.Range("F18").Value = 0.009868421052631578
$ws.Range("H18").Value = 0.1973684210526316
$ws.Range("I18").Value = 0.1085526315789474
$ws.Range("J18").Value = 0.4046052631578947
$ws.Range("K18").Value = 0.08881578947368421
$ws.Range("M18").Value = 0.03289473684210526
$ws.Range("O18").Value = 0.04605263157894737
$ws.Range("S18").Value = 0.1118421052631579
$ws.Range("F19").Value = 0.01280409731113956
$ws.Range("H19").Value = 0.2087067861715749
$ws.Range("I19").Value = 0.08386683738796415
$ws.Range("J19").Value = 0.3809218950064021
$ws.Range("K19").Value = 0.1133162612035852
$ws.Range("M19").Value = 0.01728553137003841
$ws.Range("N19").Value = 0.0006402048655569782
$ws.Range("O19").Value = 0.05633802816901409
$ws.Range("S19").Value = 0.1261203585147247
